# cs-en-us-030pct.xlsx weekly refresh: new crime data collected.
# Advances the report window by one week (1/12-1/18/2026 -> 1/19-1/25/2026,
# Number 3 -> Number 4) and refreshes the crime-complaint grid (rows 15-30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

# Replace a run of characters inside a (possibly rich-text) string cell,
# in place, preserving the remaining runs' formatting.
function Set-CharRun($range, [int]$start, [int]$length, [string]$text) {
    $chars = $range.Characters($start, $length)
    $chars.Text = $text
}

# Write a plain numeric value into a cell whose current style already has
# the correct number format (the common case: no format transition needed).
function Set-NumCell($ws, [string]$addr, $val) {
    $ws.Range($addr).Value = $val
}

# Write a numeric value into a cell that currently holds placeholder TEXT
# (shared string "0" / "***.*") and therefore needs its style switched back
# to a real numeric format. $styleSrcAddr names a nearby untouched cell that
# already carries the desired numeric style (15 = "#,##0", 14 = percent).
function Set-NumCellFixStyle($ws, [string]$addr, $val, [string]$styleSrcAddr) {
    $dst = $ws.Range($addr)
    $dst.Value = $val
    $src = $ws.Range($styleSrcAddr)
    $src.Copy()
    $dst.PasteSpecial(-4122)  # xlPasteFormats
    $ws.Application.CutCopyMode = 0
}

# Write a TEXT placeholder ("0" or "***.*") into a cell that currently holds
# a number. Excel auto-coerces a literal digit-only Value assignment back to
# a number, so route the text through a temporary formula result (which is
# always typed as text) and then freeze it to a static value; finally copy
# the neighbouring placeholder cell's style (General-formatted, right
# aligned) onto it.
function Set-TextCell($ws, [string]$addr, [string]$text) {
    $dst = $ws.Range($addr)
    $dst.Value = "x"
    $escaped = $text.Replace('"', '""')
    $dst.Formula = '="' + $escaped + '"'
    $dst.Copy()
    $dst.PasteSpecial(-4163)  # xlPasteValues - freeze formula result to a value
    $src = $ws.Range("C22")   # stable style-13 (General, right/center) source
    $src.Copy()
    $dst.PasteSpecial(-4122)  # xlPasteFormats
    $ws.Application.CutCopyMode = 0
}

# ---------------------------------------------------------------------------
# Masthead: "Volume 33   Number  3" -> "...  Number  4"
# ---------------------------------------------------------------------------
Set-CharRun $ws.Range("A8") 21 1 "4"

# ---------------------------------------------------------------------------
# "Report Covering the Week  1/12/2026  Through  1/18/2026"
#                         -> 1/19/2026            1/25/2026
# ---------------------------------------------------------------------------
Set-CharRun $ws.Range("C9") 27 9 "1/19/2026"
Set-CharRun $ws.Range("C9") 47 9 "1/25/2026"

# ---------------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------------
Set-TextCell $ws "F15" "0"
Set-TextCell $ws "G15" "0"
Set-TextCell $ws "H15" "***.*"

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
Set-NumCell $ws "C16" 1
Set-TextCell $ws "D16" "0"
Set-TextCell $ws "E16" "***.*"
Set-NumCell $ws "G16" 7
Set-NumCell $ws "H16" 0
Set-NumCell $ws "I16" 7
Set-NumCell $ws "K16" 0
Set-NumCell $ws "L16" -22.222222222222
Set-NumCell $ws "M16" -41.666666666666
Set-NumCell $ws "N16" -89.705882352941

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
Set-NumCell $ws "D17" 2
Set-NumCell $ws "E17" 0
Set-NumCell $ws "F17" 16
Set-NumCell $ws "G17" 11
Set-NumCell $ws "H17" 45.454545454545
Set-NumCell $ws "I17" 13
Set-NumCell $ws "J17" 10
Set-NumCell $ws "K17" 30
Set-NumCell $ws "L17" 8.333333333333
Set-NumCell $ws "M17" 225
Set-NumCell $ws "N17" -72.340425531914

# ---------------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------------
Set-TextCell $ws "C18" "0"
Set-NumCell $ws "E18" -100
Set-NumCell $ws "G18" 5
Set-NumCell $ws "H18" -80
Set-NumCell $ws "J18" 5
Set-NumCell $ws "K18" -80
Set-NumCell $ws "M18" -83.333333333333
Set-NumCell $ws "N18" -98.611111111111

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
Set-NumCell $ws "C19" 7
Set-NumCell $ws "D19" 5
Set-NumCell $ws "E19" 40
Set-NumCell $ws "F19" 21
Set-NumCell $ws "G19" 24
Set-NumCell $ws "H19" -12.5
Set-NumCell $ws "I19" 18
Set-NumCell $ws "J19" 22
Set-NumCell $ws "K19" -18.181818181818
Set-NumCell $ws "L19" -41.935483870967
Set-NumCell $ws "M19" 157.142857142857
Set-NumCell $ws "N19" -41.935483870967

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------------
Set-NumCellFixStyle $ws "C20" 2 "F22"
Set-TextCell $ws "D20" "0"
Set-TextCell $ws "E20" "***.*"
Set-NumCell $ws "F20" 3
Set-NumCell $ws "H20" 0
Set-NumCellFixStyle $ws "I20" 2 "F22"
Set-NumCell $ws "K20" 0
Set-NumCell $ws "L20" -66.666666666666
Set-NumCell $ws "M20" -66.666666666666
Set-NumCell $ws "N20" -92

# ---------------------------------------------------------------------------
# Row 21 - TOTAL (bold row, style unchanged throughout)
# ---------------------------------------------------------------------------
Set-NumCell $ws "C21" 12
Set-NumCell $ws "D21" 9
Set-NumCell $ws "E21" 33.333333333333
Set-NumCell $ws "F21" 48
Set-NumCell $ws "G21" 50
Set-NumCell $ws "H21" -4
Set-NumCell $ws "I21" 41
Set-NumCell $ws "J21" 46
Set-NumCell $ws "K21" -10.869565217391
Set-NumCell $ws "L21" -33.870967741935
Set-NumCell $ws "M21" 13.888888888888
Set-NumCell $ws "N21" -83.794466403162

# ---------------------------------------------------------------------------
# Row 23 - Housing
# ---------------------------------------------------------------------------
Set-NumCellFixStyle $ws "M23" -100 "H22"

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
Set-NumCell $ws "C24" 8
Set-NumCell $ws "D24" 8
Set-NumCell $ws "E24" 0
Set-NumCell $ws "I24" 37
Set-NumCell $ws "J24" 38
Set-NumCell $ws "K24" -2.631578947368
Set-NumCell $ws "L24" -21.276595744680
Set-NumCell $ws "M24" 164.285714285714

# ---------------------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------------------
Set-NumCell $ws "C25" 1
Set-NumCell $ws "D25" 2
Set-NumCell $ws "E25" -50
Set-NumCell $ws "F25" 4
Set-NumCell $ws "G25" 5
Set-NumCell $ws "H25" -20
Set-NumCell $ws "I25" 3
Set-NumCell $ws "J25" 4
Set-NumCell $ws "K25" -25

# ---------------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------------
Set-NumCell $ws "C26" 7
Set-NumCell $ws "D26" 1
Set-NumCell $ws "E26" 600
Set-NumCell $ws "F26" 28
Set-NumCell $ws "G26" 22
Set-NumCell $ws "H26" 27.272727272727
Set-NumCell $ws "I26" 26
Set-NumCell $ws "J26" 21
Set-NumCell $ws "K26" 23.809523809523
Set-NumCell $ws "L26" 62.5
Set-NumCell $ws "M26" 13.043478260869

# ---------------------------------------------------------------------------
# Row 27 - UCR Rape*
# ---------------------------------------------------------------------------
Set-TextCell $ws "F27" "0"
Set-NumCell $ws "G27" 1
Set-NumCell $ws "H27" -100

# ---------------------------------------------------------------------------
# Row 28 - Other Sex Crimes
# ---------------------------------------------------------------------------
Set-NumCellFixStyle $ws "C28" 1 "F22"
Set-NumCell $ws "F28" 4
Set-NumCell $ws "I28" 3
Set-NumCell $ws "L28" 50

# ---------------------------------------------------------------------------
# Row 29 - Shooting Vic.
# ---------------------------------------------------------------------------
Set-TextCell $ws "D29" "0"
Set-TextCell $ws "E29" "***.*"

# ---------------------------------------------------------------------------
# Row 30 - Shooting Inc.
# ---------------------------------------------------------------------------
Set-TextCell $ws "D30" "0"
Set-TextCell $ws "E30" "***.*"
